$d = $word.ActiveDocument

# Locate the paragraph that ends the DELETE section:
# "... This will delete jim's record from the table Name"
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*record from the table Name*") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph"
}

$target = $d.Range($anchor.Range.End, $anchor.Range.End)
$target.Collapse(0)

# Insert the new block of paragraphs as a raw OOXML fragment right after the
# anchor paragraph. A trailing throwaway paragraph ("–SPLITMARK–") is
# included so every paragraph above it becomes a genuine, independent
# paragraph (the very last paragraph of an InsertXML fragment always merges
# its contents into whatever paragraph used to follow the insertion point).
# That throwaway text is deleted again afterwards, which restores the
# following (pre-existing) paragraph to its original, untouched form.
$xmlFrag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">DELETE All </w:t></w:r><w:r><w:t>– to delete all records from table. But the table will not be deleted. The cols will remain intact.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">DELETE FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TableName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">DELETE * </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TableName</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p/><w:p/><w:p><w:r><w:t>–SPLITMARK–</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xmlFrag)

# Remove the throwaway split-marker text again, leaving the paragraph it was
# injected into back in its original (pre-existing) state.
$find = $d.Content
$find.Find.ClearFormatting()
$hit = $find.Find.Execute("–SPLITMARK–", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hit) {
    $find.Delete()
}

Write-Host "Inserted DELETE All block"
